# Fruta / hortaliza, semanal
# Insert a new weekly record as row 5, pushing the existing rows 5-14 down to 6-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 5:14 down to 6:15, leaving row 5 free for the new record.
$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with the new weekly price record.
$ws.Cells.Item(5, 1).Value = 10
$ws.Cells.Item(5, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(5, 3).Value = "La Araucanía"
$ws.Cells.Item(5, 4).Value = 45281
$ws.Cells.Item(5, 5).Value = 9
$ws.Cells.Item(5, 6).Value = "Fruta"
$ws.Cells.Item(5, 7).Value = 100101
$ws.Cells.Item(5, 8).Value = "Berries"
$ws.Cells.Item(5, 9).Value = 100101004
$ws.Cells.Item(5, 10).Value = "Frambuesa"
$ws.Cells.Item(5, 11).Value = "Sin especificar"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 100
$ws.Cells.Item(5, 14).Value = 3800
$ws.Cells.Item(5, 15).Value = 3800
$ws.Cells.Item(5, 16).Value = 3800
$ws.Cells.Item(5, 17).Value = "$/envase 1 kilo"
$ws.Cells.Item(5, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(5, 19).Value = 3800
$ws.Cells.Item(5, 20).Value = 1
